# Auto-generated edit script: refresh scraped schedule data for Línea 141
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 'Última actualización: 07:48:35'

$ws.Range("A3").Value = 'Total filas: 81'

$ws.Range("A46").Value = '06:33:46'
$ws.Range("B46").Value = '07:59'
$ws.Range("C46").Value = '11_ETCHEVERRY'
$ws.Range("D46").Value = 86
$ws.Range("E46").Value = 'LP1912'

$ws.Range("A47").Value = '07:12:53'
$ws.Range("B47").Value = '07:59'
$ws.Range("C47").Value = '23_HERNANDEZ'
$ws.Range("D47").Value = 47
$ws.Range("E47").Value = 'LP1912'

$ws.Range("A57").Value = '07:48:35'
$ws.Range("B57").Value = '08:14'
$ws.Range("C57").Value = '10_OLMOS'
$ws.Range("D57").Value = 26
$ws.Range("E57").Value = 'LP1912'

$ws.Range("A58").Value = '07:36:59'
$ws.Range("B58").Value = '08:14'
$ws.Range("C58").Value = '17_ROMERO'
$ws.Range("D58").Value = 38
$ws.Range("E58").Value = 'LP1912'

$ws.Range("A59").Value = '06:16:15'
$ws.Range("B59").Value = '08:15'
$ws.Range("C59").Value = '17_ROMERO'
$ws.Range("D59").Value = 119
$ws.Range("E59").Value = 'LP1912'

$ws.Range("A60").Value = '07:36:59'
$ws.Range("B60").Value = '08:25'
$ws.Range("C60").Value = '15X38_ABASTO'
$ws.Range("D60").Value = 49
$ws.Range("E60").Value = 'LP1912'

$ws.Range("A61").Value = '06:33:46'
$ws.Range("B61").Value = '08:26'
$ws.Range("C61").Value = '15X38_ABASTO'
$ws.Range("D61").Value = 113
$ws.Range("E61").Value = 'LP1912'

$ws.Range("A62").Value = '06:33:46'
$ws.Range("B62").Value = '08:27'
$ws.Range("C62").Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Range("D62").Value = 114
$ws.Range("E62").Value = 'LP1912'

$ws.Range("A63").Value = '06:45:50'
$ws.Range("B63").Value = '08:29'
$ws.Range("C63").Value = '14_ABASTO'
$ws.Range("D63").Value = 104
$ws.Range("E63").Value = 'LP1912'

$ws.Range("A64").Value = '07:36:59'
$ws.Range("B64").Value = '08:30'
$ws.Range("C64").Value = '16_P MOR-SANTA ANA'
$ws.Range("D64").Value = 54
$ws.Range("E64").Value = 'LP1912'

$ws.Range("A65").Value = '06:33:46'
$ws.Range("B65").Value = '08:31'
$ws.Range("C65").Value = '16_P MOR-SANTA ANA'
$ws.Range("D65").Value = 118
$ws.Range("E65").Value = 'LP1912'

$ws.Range("A66").Value = '06:45:50'
$ws.Range("B66").Value = '08:38'
$ws.Range("C66").Value = '215C_EL PATO'
$ws.Range("D66").Value = 113
$ws.Range("E66").Value = 'LP1912'

$ws.Range("A67").Value = '07:48:35'
$ws.Range("B67").Value = '08:39'
$ws.Range("C67").Value = '215C_EL PATO'
$ws.Range("D67").Value = 51
$ws.Range("E67").Value = 'LP1912'

$ws.Range("A68").Value = '07:12:53'
$ws.Range("B68").Value = '08:43'
$ws.Range("C68").Value = '10_OLMOS'
$ws.Range("D68").Value = 91
$ws.Range("E68").Value = 'LP1912'

$ws.Range("A69").Value = '07:48:35'
$ws.Range("B69").Value = '08:44'
$ws.Range("C69").Value = '10_OLMOS'
$ws.Range("D69").Value = 56
$ws.Range("E69").Value = 'LP1912'

$ws.Range("A70").Value = '07:12:53'
$ws.Range("B70").Value = '08:49'
$ws.Range("C70").Value = '215A_EL PATO'
$ws.Range("D70").Value = 97
$ws.Range("E70").Value = 'LP1912'

$ws.Range("A71").Value = '07:48:35'
$ws.Range("B71").Value = '08:50'
$ws.Range("C71").Value = '215A_EL PATO'
$ws.Range("D71").Value = 62
$ws.Range("E71").Value = 'LP1912'

$ws.Range("A72").Value = '07:12:53'
$ws.Range("B72").Value = '08:59'
$ws.Range("C72").Value = '215B_EL PATO'
$ws.Range("D72").Value = 107
$ws.Range("E72").Value = 'LP1912'

$ws.Range("A73").Value = '07:36:59'
$ws.Range("B73").Value = '09:01'
$ws.Range("C73").Value = '17X38_ROMERO'
$ws.Range("D73").Value = 85
$ws.Range("E73").Value = 'LP1912'

$ws.Range("A74").Value = '07:36:59'
$ws.Range("B74").Value = '09:02'
$ws.Range("C74").Value = '23_HERNANDEZ'
$ws.Range("D74").Value = 86
$ws.Range("E74").Value = 'LP1912'

$ws.Range("A75").Value = '07:12:53'
$ws.Range("B75").Value = '09:02'
$ws.Range("C75").Value = '17X38_ROMERO'
$ws.Range("D75").Value = 110
$ws.Range("E75").Value = 'LP1912'

$ws.Range("A76").Value = '07:36:59'
$ws.Range("B76").Value = '09:04'
$ws.Range("C76").Value = '16_SANTA ANA'
$ws.Range("D76").Value = 88
$ws.Range("E76").Value = 'LP1912'

$ws.Range("A77").Value = '07:48:35'
$ws.Range("B77").Value = '09:08'
$ws.Range("C77").Value = '16_SANTA ANA'
$ws.Range("D77").Value = 80
$ws.Range("E77").Value = 'LP1912'

$ws.Range("A78").Value = '07:36:59'
$ws.Range("B78").Value = '09:14'
$ws.Range("C78").Value = '11_ETCHEVERRY'
$ws.Range("D78").Value = 98
$ws.Range("E78").Value = 'LP1912'

$ws.Range("A79").Value = '07:36:59'
$ws.Range("B79").Value = '09:14'
$ws.Range("C79").Value = '15_ABASTO'
$ws.Range("D79").Value = 98
$ws.Range("E79").Value = 'LP1912'

$ws.Range("A80").Value = '07:36:59'
$ws.Range("B80").Value = '09:16'
$ws.Range("C80").Value = '27_EL RETIRO'
$ws.Range("D80").Value = 100
$ws.Range("E80").Value = 'LP1912'

$ws.Range("A81").Value = '07:36:59'
$ws.Range("B81").Value = '09:26'
$ws.Range("C81").Value = '215_EL PELIGRO'
$ws.Range("D81").Value = 110
$ws.Range("E81").Value = 'LP1912'

$ws.Range("A82").Value = '07:48:35'
$ws.Range("B82").Value = '09:27'
$ws.Range("C82").Value = '215_EL PELIGRO'
$ws.Range("D82").Value = 99
$ws.Range("E82").Value = 'LP1912'

$ws.Range("A83").Value = '07:36:59'
$ws.Range("B83").Value = '09:30'
$ws.Range("C83").Value = '16_P MOR-SANTA ANA'
$ws.Range("D83").Value = 114
$ws.Range("E83").Value = 'LP1912'

$ws.Range("A84").Value = '07:48:35'
$ws.Range("B84").Value = '09:31'
$ws.Range("C84").Value = '16_P MOR-SANTA ANA'
$ws.Range("D84").Value = 103
$ws.Range("E84").Value = 'LP1912'

$ws.Range("A85").Value = '07:48:35'
$ws.Range("B85").Value = '09:39'
$ws.Range("C85").Value = '15_ABASTO'
$ws.Range("D85").Value = 111
$ws.Range("E85").Value = 'LP1912'

$ws.Range("A86").Value = '07:48:35'
$ws.Range("B86").Value = '09:44'
$ws.Range("C86").Value = '14_ABASTO'
$ws.Range("D86").Value = 116
$ws.Range("E86").Value = 'LP1912'

# ---- Sheet 2: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)

$ws.Range("A2").Value = 'Última actualización: 07:48:35'

$ws.Range("A3").Value = 'Total filas: 15'

$ws.Range("A15").Value = '07:48:35'
$ws.Range("B15").Value = '08:39'
$ws.Range("C15").Value = '215C_EL PATO'
$ws.Range("D15").Value = 51
$ws.Range("E15").Value = 'LP1912'

$ws.Range("A16").Value = '07:12:53'
$ws.Range("B16").Value = '08:49'
$ws.Range("C16").Value = '215A_EL PATO'
$ws.Range("D16").Value = 97
$ws.Range("E16").Value = 'LP1912'

$ws.Range("A17").Value = '07:48:35'
$ws.Range("B17").Value = '08:50'
$ws.Range("C17").Value = '215A_EL PATO'
$ws.Range("D17").Value = 62
$ws.Range("E17").Value = 'LP1912'

$ws.Range("A18").Value = '07:12:53'
$ws.Range("B18").Value = '08:59'
$ws.Range("C18").Value = '215B_EL PATO'
$ws.Range("D18").Value = 107
$ws.Range("E18").Value = 'LP1912'

$ws.Range("A19").Value = '07:36:59'
$ws.Range("B19").Value = '09:26'
$ws.Range("C19").Value = '215_EL PELIGRO'
$ws.Range("D19").Value = 110
$ws.Range("E19").Value = 'LP1912'

$ws.Range("A20").Value = '07:48:35'
$ws.Range("B20").Value = '09:27'
$ws.Range("C20").Value = '215_EL PELIGRO'
$ws.Range("D20").Value = 99
$ws.Range("E20").Value = 'LP1912'

# ---- Sheet 3: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)

$ws.Range("A2").Value = 'Última actualización: 07:48:35'

$ws.Range("A3").Value = 'Total filas: 8'

$ws.Range("A10").Value = '07:48:35'
$ws.Range("B10").Value = '08:25'
$ws.Range("C10").Value = '215C_LA PLATA'
$ws.Range("D10").Value = 37
$ws.Range("E10").Value = 'L6203'

$ws.Range("A11").Value = '07:36:59'
$ws.Range("B11").Value = '08:27'
$ws.Range("C11").Value = '215C_LA PLATA'
$ws.Range("D11").Value = 51
$ws.Range("E11").Value = 'L6203'

$ws.Range("A12").Value = '07:36:59'
$ws.Range("B12").Value = '08:51'
$ws.Range("C12").Value = '215A_LA PLATA'
$ws.Range("D12").Value = 75
$ws.Range("E12").Value = 'L6173'

$ws.Range("A13").Value = '07:48:35'
$ws.Range("B13").Value = '08:52'
$ws.Range("C13").Value = '215A_LA PLATA'
$ws.Range("D13").Value = 64
$ws.Range("E13").Value = 'L6173'

